$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new values (custom accuracy run)
# Row 2
$ws.Range("A2").Value = 45037.50694444445
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0

# Row 3
$ws.Range("A3").Value = 45037.51388888889
$ws.Range("B3").Value = 0.468
$ws.Range("C3").Value = 0.261
$ws.Range("D3").Value = 0.174
$ws.Range("E3").Value = 0.365
$ws.Range("F3").Value = 0.064
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.637
$ws.Range("I3").Value = 0.218
$ws.Range("J3").Value = 0.229
$ws.Range("K3").Value = 0.264
$ws.Range("L3").Value = 0.374
$ws.Range("M3").Value = 0.347
$ws.Range("N3").Value = 0.286
$ws.Range("O3").Value = 0.248
$ws.Range("P3").Value = 0.471
$ws.Range("Q3").Value = 0.272
$ws.Range("R3").Value = 0.238
$ws.Range("S3").Value = 0.202
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0.525
$ws.Range("V3").Value = 0.271
$ws.Range("W3").Value = 0.74
$ws.Range("X3").Value = 0.248
$ws.Range("Y3").Value = 0.124
$ws.Range("Z3").Value = 0.265
$ws.Range("AA3").Value = 0.32
$ws.Range("AB3").Value = 1.002
$ws.Range("AC3").Value = 0.213
$ws.Range("AD3").Value = 0.605
$ws.Range("AE3").Value = 0.621
$ws.Range("AF3").Value = 0.472
$ws.Range("AG3").Value = 0.264
$ws.Range("AH3").Value = 0.25

# Row 4
$ws.Range("A4").Value = 45037.52083333334
$ws.Range("B4").Value = 23.053
$ws.Range("C4").Value = 17.232
$ws.Range("D4").Value = 0.912
$ws.Range("E4").Value = 49.648
$ws.Range("F4").Value = 40.901
$ws.Range("G4").Value = 18.348
$ws.Range("H4").Value = 60.958
$ws.Range("I4").Value = 27.652
$ws.Range("J4").Value = 12.513
$ws.Range("K4").Value = 18.644
$ws.Range("L4").Value = 20.115
$ws.Range("M4").Value = 21.36
$ws.Range("N4").Value = 5.881
$ws.Range("O4").Value = 17.963
$ws.Range("P4").Value = 25.788
$ws.Range("Q4").Value = 14.914
$ws.Range("R4").Value = 0.391
$ws.Range("S4").Value = 0.827
$ws.Range("T4").Value = 267.819
$ws.Range("U4").Value = 49.937
$ws.Range("V4").Value = 16.593
$ws.Range("W4").Value = 34.202
$ws.Range("X4").Value = 17.818
$ws.Range("Y4").Value = 2.405
$ws.Range("Z4").Value = 31.292
$ws.Range("AA4").Value = 14.721
$ws.Range("AB4").Value = 13.398
$ws.Range("AC4").Value = 15.149
$ws.Range("AD4").Value = 21.266
$ws.Range("AE4").Value = 0.453
$ws.Range("AF4").Value = 55.201
$ws.Range("AG4").Value = 9.585000000000001
$ws.Range("AH4").Value = 20.682

# Row 5
$ws.Range("A5").Value = 45037.52777777778
$ws.Range("B5").Value = 10.08
$ws.Range("C5").Value = 7.49
$ws.Range("D5").Value = 0.46
$ws.Range("E5").Value = 21.53
$ws.Range("F5").Value = 17.7
$ws.Range("G5").Value = 8.15
$ws.Range("H5").Value = 34.1
$ws.Range("I5").Value = 12
$ws.Range("J5").Value = 5.49
$ws.Range("K5").Value = 8.130000000000001
$ws.Range("L5").Value = 8.800000000000001
$ws.Range("M5").Value = 9.4
$ws.Range("N5").Value = 2.59
$ws.Range("O5").Value = 7.83
$ws.Range("P5").Value = 11.34
$ws.Range("Q5").Value = 6.52
$ws.Range("R5").Value = 0.25
$ws.Range("S5").Value = 0.42
$ws.Range("T5").Value = 113.15
$ws.Range("U5").Value = 21.96
$ws.Range("V5").Value = 7.23
$ws.Range("W5").Value = 15.18
$ws.Range("X5").Value = 7.77
$ws.Range("Y5").Value = 1.06
$ws.Range("Z5").Value = 16.28
$ws.Range("AA5").Value = 6.44
$ws.Range("AB5").Value = 6.02
$ws.Range("AC5").Value = 6.63
$ws.Range("AD5").Value = 9.33
$ws.Range("AE5").Value = 0.36
$ws.Range("AF5").Value = 31.33
$ws.Range("AG5").Value = 4.24
$ws.Range("AH5").Value = 8.99

# Remove row 6 (dataset now has one fewer sample row)
$ws.Rows.Item(6).Delete()

# Widen several data columns by one character unit (custom accuracy formatting)
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(31).ColumnWidth = 6.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
